$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 3 updates (existing "App Dev" task) ---
# D3: "[()]" -> "[]"
$ws.Range("D3").Value = "[]"
# E3: "[()]" -> new activity repr string
$ws.Range("E3").Value = "[(<object.Activity.Activity object at 0x7f5311b13ef0>, 'FS', 0)]"
# F3: baseline start serial changes
$ws.Range("F3").Value = 42097.6407448236
# G3: text "0" (shared string, not numeric) - force text without leaving a styled cell
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "0"
$ws.Range("G3").ClearFormats()
# H3: baseline end serial changes
$ws.Range("H3").Value = 42102.6407448236
# I3: text "0" (shared string)
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "0"
$ws.Range("I3").ClearFormats()

# --- Row 4: new task "Testing" ---
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Testing"
$ws.Range("C4").Value = "(1, 1)"
$ws.Range("D4").Value = "[(<object.Activity.Activity object at 0x7f5311b13eb8>, 'FS', 0)]"
$ws.Range("E4").Value = "[]"
$ws.Range("F4").Value = 42097.6407448238
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "0"
$ws.Range("G4").ClearFormats()
$ws.Range("H4").Value = 42107.6407448238
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "0"
$ws.Range("I4").ClearFormats()
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("M4").Value = 10000
$ws.Range("N4").Value = 0
